# Generate Report for Handback
# - Marks the Overview sheet rows for a.md.md / b.md.md as handed back.
# - Records the handback target/file + datetime for each locale sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet: status text changes for a.md.md / b.md.md in both
# locale columns (zh-cn = B, de-de = C).
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = "Handed back: in sync with en-US"
$overview.Range("C2").Value = "Handed back: in sync with en-US"
$overview.Range("B3").Value = "Handed back: in sync with en-US"
$overview.Range("C3").Value = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# zh-cn sheet: populate "Latest Target File" (E) / "Latest Handback
# File" (F) hyperlinks and the "Latest Handback DateTime" (G) for the
# two localized rows.
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

$zhcn.Hyperlinks.Add($zhcn.Range("E2"), "https://github.com/OpenLocalizationTest/oltest/blob/aa0916ab69ae653efbc511c6b1d264e9bfb36354/e2e/a.md.md", "", "", "a.md.md")
$zhcn.Hyperlinks.Add($zhcn.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/aa0916ab69ae653efbc511c6b1d264e9bfb36354/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/a.md.370104d57010292b5765347db07350cde3a977e8.zh-cn.xlf", "", "", "a.md.370104d57010292b5765347db07350cde3a977e8.zh-cn.xlf")
$zhcn.Range("G2").Value = "2016-01-25 07:22:11"

$zhcn.Hyperlinks.Add($zhcn.Range("E3"), "https://github.com/OpenLocalizationTest/oltest/blob/aa0916ab69ae653efbc511c6b1d264e9bfb36354/e2e/a.md.md", "", "", "a.md.md")
$zhcn.Hyperlinks.Add($zhcn.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/aa0916ab69ae653efbc511c6b1d264e9bfb36354/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/a.md.370104d57010292b5765347db07350cde3a977e8.zh-cn.xlf", "", "", "a.md.370104d57010292b5765347db07350cde3a977e8.zh-cn.xlf")
$zhcn.Range("G3").Value = "2016-01-25 07:22:11"

# ---------------------------------------------------------------------
# de-de sheet: same pattern as zh-cn, with de-de specific file names and
# the de-de handback datetime.
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

$dede.Hyperlinks.Add($dede.Range("E2"), "https://github.com/OpenLocalizationTest/oltest/blob/aa0916ab69ae653efbc511c6b1d264e9bfb36354/e2e/a.md.md", "", "", "a.md.md")
$dede.Hyperlinks.Add($dede.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/ec52321ed6d6cc86f8419df615acb8f223751d06/ol-handback/OpenLocalizationTestOrg/oltest.de-de/yuwzho/a.md.370104d57010292b5765347db07350cde3a977e8.de-de.xlf", "", "", "a.md.370104d57010292b5765347db07350cde3a977e8.de-de.xlf")
$dede.Range("G2").Value = "2016-01-25 07:22:32"

$dede.Hyperlinks.Add($dede.Range("E3"), "https://github.com/OpenLocalizationTest/oltest/blob/aa0916ab69ae653efbc511c6b1d264e9bfb36354/e2e/a.md.md", "", "", "a.md.md")
$dede.Hyperlinks.Add($dede.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/ec52321ed6d6cc86f8419df615acb8f223751d06/ol-handback/OpenLocalizationTestOrg/oltest.de-de/yuwzho/a.md.370104d57010292b5765347db07350cde3a977e8.de-de.xlf", "", "", "a.md.370104d57010292b5765347db07350cde3a977e8.de-de.xlf")
$dede.Range("G3").Value = "2016-01-25 07:22:32"
